$wb = $excel.ActiveWorkbook

# --- Data sheet: title text + 2015 -> 2016 fixes, row height, active tab/selection ---
$wsData = $wb.Worksheets.Item("Data")

# Fix the report title (shared string) - "September 2015" -> "September 2016"
$wsData.Range("A1").Value = "The proportion of early childhood education and care services in Australia with a quality rating, by quality rating level and jurisdiction, March 2013 and September 2016"

# Fix the year label in the second data block - 2015 -> 2016
$wsData.Range("A9").Value = 2016

# Row 9 height correction (15.65 -> 15)
$wsData.Rows.Item(9).RowHeight = 15

# Make "Data" the active/selected sheet with A1 selected
$wsData.Activate()
[void]$wsData.Range("A1").Select()

# --- Description sheet: fix corrupted font on B2 (reuse the correct Arial font from B1) ---
$wsDesc = $wb.Worksheets.Item("Description")
$wsDesc.Range("B1").Copy()
[void]$wsDesc.Range("B2").PasteSpecial(-4122)
